$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture all source values we need BEFORE mutating the sheet ---
$q_solidColor   = $ws.Cells.Item(4,3).Value()   # "Problem:...solid color..."
$q_cutInOut     = $ws.Cells.Item(8,3).Value()   # "Problem:...cut in & out..."
$q_static       = $ws.Cells.Item(10,3).Value()  # "Problem:...static, wavy lines..."

$possibleProblem = $ws.Cells.Item(6,2).Value()  # "Possible_Problem"
$pp55            = $ws.Cells.Item(6,3).Value()  # "Possible_Problem:55%Failed Camera..."

$yes = $ws.Cells.Item(7,2).Value()              # "Yes"
$no  = $ws.Cells.Item(8,2).Value()              # "No"

$pp70Failed      = $ws.Cells.Item(7,3).Value()  # "Possible_Problem:70% Failed Camera..."
$pp60Audio       = $ws.Cells.Item(8,3).Value()  # cut in & out question text (string idx 16) - NOT USED directly
$pp50Audio       = $ws.Cells.Item(9,3).Value()  # "Possible_Problem:50% Audio Control Unit (Radio)..."
$pp60Cable       = $ws.Cells.Item(10,3).Value() # static question text (string idx 18) - NOT USED directly
$pp60CableConn   = $ws.Cells.Item(11,3).Value() # "Possible_Problem:60% Camera Cable / Connection..."
$pp70Failed25    = $ws.Cells.Item(12,3).Value() # "Possible_Problem:70%Failed Camera\n25%..."

# Also capture the two "branch" question texts used as A-values for the new rows
$q_cutInOutText = $ws.Cells.Item(8,3).Value()   # row8,col C = string idx 16
$q_staticText   = $ws.Cells.Item(10,3).Value()  # row10,col C = string idx 18

# --- Insert 3 new rows starting at row 9 (shifts old rows 9-12 down to 12-15) ---
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

# --- Row 7: rewrite in place (was Yes/70% branch, becomes Possible_Problem/55% branch) ---
$ws.Cells.Item(7,1).Value = $q_solidColor
$ws.Cells.Item(7,2).Value = $possibleProblem
$ws.Cells.Item(7,3).Value = $pp55
$ws.Rows.Item(7).RowHeight = 172.8

# --- Row 8: rewrite in place (was No/60% branch, becomes Yes/70% branch) ---
$ws.Cells.Item(8,1).Value = $q_solidColor
$ws.Cells.Item(8,2).Value = $yes
$ws.Cells.Item(8,3).Value = $pp70Failed
$ws.Rows.Item(8).RowHeight = 216

# --- Row 9 (new): solid-color / No -> cut-in-&-out question ---
$ws.Cells.Item(9,1).Value = $q_solidColor
$ws.Cells.Item(9,2).Value = $no
$ws.Cells.Item(9,3).Value = $q_cutInOutText
$ws.Rows.Item(9).RowHeight = 216

# --- Row 10 (new): cut-in-&-out / Possible_Problem -> 55% text ---
$ws.Cells.Item(10,1).Value = $q_cutInOutText
$ws.Cells.Item(10,2).Value = $possibleProblem
$ws.Cells.Item(10,3).Value = $pp55
$ws.Rows.Item(10).RowHeight = 172.8

# --- Row 11 (new): cut-in-&-out / Yes -> 50% Audio text ---
$ws.Cells.Item(11,1).Value = $q_cutInOutText
$ws.Cells.Item(11,2).Value = $yes
$ws.Cells.Item(11,3).Value = $pp50Audio
$ws.Rows.Item(11).RowHeight = 216

# --- Row 12 (new): cut-in-&-out / No -> static question ---
$ws.Cells.Item(12,1).Value = $q_cutInOutText
$ws.Cells.Item(12,2).Value = $no
$ws.Cells.Item(12,3).Value = $q_staticText
$ws.Rows.Item(12).RowHeight = 230.4

# --- Row 13 (new): static / Possible_Problem -> 55% text ---
$ws.Cells.Item(13,1).Value = $q_staticText
$ws.Cells.Item(13,2).Value = $possibleProblem
$ws.Cells.Item(13,3).Value = $pp55
$ws.Rows.Item(13).RowHeight = 172.8

# --- Row 14 (new): static / Yes -> 60% Camera Cable text ---
$ws.Cells.Item(14,1).Value = $q_staticText
$ws.Cells.Item(14,2).Value = $yes
$ws.Cells.Item(14,3).Value = $pp60CableConn
$ws.Rows.Item(14).RowHeight = 216

# --- Row 15 (new): static / No -> 70% Failed (25%) text ---
$ws.Cells.Item(15,1).Value = $q_staticText
$ws.Cells.Item(15,2).Value = $no
$ws.Cells.Item(15,3).Value = $pp70Failed25
$ws.Rows.Item(15).RowHeight = 172.8

# --- Apply wrap-text style (style index 1 in the original) to column C for the new rows ---
$ws.Range("C7:C15").WrapText = $true

# --- Update sheet view (scrolled down, C7 selected) ---
$ws.Range("C7").Select()
$ws.Application.ActiveWindow.ScrollRow = 6
